$d = $word.ActiveDocument
$t = $d.Tables(1)

$t.Cell(2,2).Range.Text = "101 (30.7)"
$t.Cell(3,2).Range.Text = "102 (31.0)"
$t.Cell(4,2).Range.Text = "98 (29.8)"
$t.Cell(5,2).Range.Text = "23 (7.0)"
$t.Cell(6,2).Range.Text = "2 (0.6)"

$beforeRow = $t.Rows(8)
$newRow = $t.Rows.Add($beforeRow)
$t.Cell(8,2).Range.Text = "2 (0.6)"

$t.Cell(9,2).Range.Text = "329 (100.0)"
